$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply scraped cryptos-list update (prices / 1h-volume %, plus a rank swap at rows 31-32).
# NumberFormat "@" forces literal text entry (so values like "1.005" or "342.49" are not
# auto-parsed as numbers by Excel), then Style is reset to "Normal" so no stray number-format
# style sticks to the cell (matches the source inline-string cells, which carry no style).
function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "30.271.03"
Set-TextCell "E2" "  +1.80%  "
Set-TextCell "D3" "2.087.72"
Set-TextCell "E3" "  -0.56%  "
Set-TextCell "D4" "1.005"
Set-TextCell "E4" "  -0.25%  "
Set-TextCell "D5" "342.49"
Set-TextCell "E5" "  -0.75%  "
Set-TextCell "E7" "  +1.61%  "
Set-TextCell "D8" "0.4400"
Set-TextCell "E8" "  -0.40%  "
Set-TextCell "D9" "54.42"
Set-TextCell "E9" "  +3.20%  "
Set-TextCell "D10" "0.09340"
Set-TextCell "E10" "  +1.33%  "
Set-TextCell "D11" "1.166"
Set-TextCell "E11" "  -0.57%  "
Set-TextCell "D12" "24.73"
Set-TextCell "E12" "  -1.01%  "
Set-TextCell "D13" "8.586"
Set-TextCell "E13" "  +3.70%  "
Set-TextCell "D14" "6.886"
Set-TextCell "E14" "  +1.82%  "
Set-TextCell "D15" "2.049.02"
Set-TextCell "E15" "  -2.59%  "
Set-TextCell "D16" "101.35"
Set-TextCell "E16" "  +1.65%  "
Set-TextCell "E17" "  +0.26%  "
Set-TextCell "D18" "1.005"
Set-TextCell "E18" "  -0.34%  "
Set-TextCell "D19" "21.09"
Set-TextCell "E19" "  +1.02%  "
Set-TextCell "D20" "0.06670"
Set-TextCell "E20" "  +0.58%  "
Set-TextCell "D21" "6.317"
Set-TextCell "E21" "  +1.85%  "
Set-TextCell "E22" "  -0.37%  "
Set-TextCell "D23" "30.285.00"
Set-TextCell "E23" "  +1.71%  "
Set-TextCell "D24" "12.48"
Set-TextCell "E24" "  -1.18%  "
Set-TextCell "D25" "2.315"
Set-TextCell "E25" "  -0.09%  "
Set-TextCell "D26" "21.77"
Set-TextCell "E26" "  -0.68%  "
Set-TextCell "D27" "162.09"
Set-TextCell "E27" "  +0.04%  "
Set-TextCell "D28" "2.505"
Set-TextCell "E28" "  -1.17%  "
Set-TextCell "D29" "132.98"
Set-TextCell "E29" "  +0.00%  "
Set-TextCell "D30" "1.127"
Set-TextCell "E30" "  -0.46%  "
Set-TextCell "B31" "Stellar"
Set-TextCell "C31" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D31" "0.1045"
Set-TextCell "E31" "  -0.51%  "
Set-TextCell "B32" "ARBITRUM"
Set-TextCell "C32" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D32" "1.653"
Set-TextCell "E32" "  -0.37%  "
Set-TextCell "D33" "6.213"
Set-TextCell "E33" "  +0.50%  "
Set-TextCell "D34" "6.652"
Set-TextCell "E34" "  +10.01%  "
Set-TextCell "D35" "3.872"
Set-TextCell "E35" "  -1.64%  "
Set-TextCell "E36" "  -2.69%  "
Set-TextCell "D37" "0.02627"
Set-TextCell "E37" "  +2.19%  "
Set-TextCell "D38" "0.06784"
Set-TextCell "E38" "  +0.64%  "
Set-TextCell "D39" "0.6960"
Set-TextCell "E39" "  +1.09%  "
Set-TextCell "D40" "1.340"
Set-TextCell "E40" "  +2.70%  "
Set-TextCell "D41" "12.49"
Set-TextCell "E41" "  +0.18%  "
Set-TextCell "D42" "0.2203"
Set-TextCell "E42" "  -1.45%  "
Set-TextCell "D43" "0.6788"
Set-TextCell "E43" "  +1.85%  "
Set-TextCell "D44" "14.25"
Set-TextCell "E44" "  +0.07%  "
Set-TextCell "D45" "2.327"
Set-TextCell "E45" "  +0.72%  "
Set-TextCell "E46" "  -0.29%  "
Set-TextCell "E47" "  +18.03%  "
Set-TextCell "D48" "3.629"
Set-TextCell "E48" "  +0.43%  "
Set-TextCell "D49" "0.00000000351"
Set-TextCell "E49" "  +1.53%  "
Set-TextCell "D50" "1.212"
Set-TextCell "E50" "  +8.18%  "
Set-TextCell "D51" "1.214"
Set-TextCell "E51" "  -0.65%  "
